$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row duplicating the existing credentials (row 2) as row 3
$ws.Range("A3").Value = "Admin"
$ws.Range("B3").Value = "admin123"

# Update the selected cell in the sheet view
$ws.Range("E7").Select()
